$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "mean" row (row 3)
$ws.Range("B3").Value = 4041.924321669676
$ws.Range("D3").Value = 271.2254994103271

# Update the "std" row (row 4)
$ws.Range("B4").Value = 1819.763485645279
$ws.Range("D4").Value = 233.8897319487942

# Update the "min" row (row 5)
$ws.Range("B5").Value = 146.014

# Update the "25%" row (row 6)
$ws.Range("B6").Value = 2604.118999999999
$ws.Range("D6").Value = 200.001

# Update the "50%" row (row 7)
$ws.Range("B7").Value = 3643.011500000001
$ws.Range("D7").Value = 280.001

# Update the "75%" row (row 8)
$ws.Range("B8").Value = 5494.01
$ws.Range("D8").Value = 295.004

# Update the "max" row (row 9)
$ws.Range("B9").Value = 13187.064
$ws.Range("D9").Value = 2180

# Update the Total / Residential / Community / IGA sums and percentages (rows 10-13)
$ws.Range("F10").Value = 2124435423.468

$ws.Range("G11").Value = 0.8063540400844776

$ws.Range("F12").Value = 142556122.49
$ws.Range("G12").Value = 0.06710306226078952

$ws.Range("G13").Value = 0.1265428976547328
